$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that can look like plain numbers
# (e.g. "596.20"). Excel auto-converts such text to a numeric Value, which
# would drop the trailing zero / exact decimal formatting. Force text type
# via NumberFormat "@" while writing, then restore the default "Normal"
# style so the cell itself is left without an explicit style, matching the
# original workbook.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '62.486.11'
$ws.Range('E2').Value = '  +1.37%  '
Set-TextValue 'D3' '3.019.81'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '596.20'
$ws.Range('E5').Value = '  +1.35%  '
Set-TextValue 'D6' '150.51'
$ws.Range('E6').Value = '  +6.14%  '
$ws.Range('E7').Value = '  -0.04%  '
Set-TextValue 'D8' '3.017.83'
$ws.Range('E8').Value = '  +1.38%  '
Set-TextValue 'D9' '0.518'
$ws.Range('E9').Value = '  -0.26%  '
Set-TextValue 'D10' '6.38'
$ws.Range('E10').Value = '  +10.58%  '
$ws.Range('E11').Value = '  +4.72%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  +3.07%  '
Set-TextValue 'D14' '34.58'
$ws.Range('E14').Value = '  +1.74%  '
Set-TextValue 'D16' '3.519.57'
$ws.Range('E16').Value = '  +1.47%  '
Set-TextValue 'D17' '62.479.34'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('E18').Value = '  -0.01%  '
Set-TextValue 'D19' '3.020.06'
$ws.Range('E19').Value = '  +1.47%  '
Set-TextValue 'D20' '449.09'
$ws.Range('E20').Value = '  -0.19%  '
Set-TextValue 'D21' '14.17'
$ws.Range('E21').Value = '  +2.10%  '
Set-TextValue 'D22' '0.689'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('E23').Value = '  +1.84%  '
Set-TextValue 'D24' '82.36'
$ws.Range('E24').Value = '  +1.46%  '
Set-TextValue 'D25' '10.88'
$ws.Range('E25').Value = '  +10.98%  '
$ws.Range('E26').Value = '  +4.23%  '
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +3.02%  '
$ws.Range('E30').Value = '  +6.84%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('E34').Value = '  +2.73%  '
Set-TextValue 'D35' '0.0₃0853'
$ws.Range('E35').Value = '  +10.20%  '
$ws.Range('E36').Value = '  +1.41%  '
Set-TextValue 'D37' '5.84'
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('E38').Value = '  +8.91%  '
$ws.Range('E39').Value = '  -0.22%  '
Set-TextValue 'D40' '50.10'
$ws.Range('E40').Value = '  -0.04%  '
Set-TextValue 'D41' '9.04'
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('E42').Value = '  +4.56%  '
$ws.Range('E43').Value = '  +9.27%  '
Set-TextValue 'D44' '40.26'
$ws.Range('E44').Value = '  +8.99%  '
Set-TextValue 'D45' '390.36'
$ws.Range('E45').Value = '  +0.70%  '
Set-TextValue 'D46' '0.0355'
$ws.Range('E46').Value = '  +0.30%  '
Set-TextValue 'D47' '2.735.36'
$ws.Range('E47').Value = '  +0.36%  '
Set-TextValue 'D48' '132.81'
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('E51').Value = '  -0.11%  '
